$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '52.406.34'
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.913.38'
$ws.Range('E3').Value = '  +3.74%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '352.78'
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '112.46'
$ws.Range('E6').Value = '  +0.93%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.638'
$ws.Range('E9').Value = '  +0.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.99'
$ws.Range('E10').Value = '  -1.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0868'
$ws.Range('E11').Value = '  +3.14%  '
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.87'
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.83'
$ws.Range('E14').Value = '  +0.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.370.16'
$ws.Range('E15').Value = '  +3.64%  '
$ws.Range('E16').Value = '  +6.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.921.59'
$ws.Range('E17').Value = '  +3.89%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '52.404.63'
$ws.Range('E18').Value = '  +1.13%  '
$ws.Range('E19').Value = '  -0.40%  '
$ws.Range('E20').Value = '  +3.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.24'
$ws.Range('E21').Value = '  +4.20%  '
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.03'
$ws.Range('E23').Value = '  +0.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '270.29'
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.78'
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.79'
$ws.Range('E26').Value = '  +2.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.170'
$ws.Range('E27').Value = '  +4.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.67'
$ws.Range('E29').Value = '  +2.63%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.66'
$ws.Range('E30').Value = '  +8.41%  '
$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.37'
$ws.Range('E31').Value = '  +12.93%  '
$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '37.95'
$ws.Range('E32').Value = '  -1.59%  '
$ws.Range('E33').Value = '  +0.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0985'
$ws.Range('E34').Value = '  +11.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '53.42'
$ws.Range('E35').Value = '  +1.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0451'
$ws.Range('E36').Value = '  +1.74%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.33'
$ws.Range('E38').Value = '  +5.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.84'
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.87'
$ws.Range('E40').Value = '  +14.64%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.07'
$ws.Range('E41').Value = '  +2.90%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '23.69'
$ws.Range('E42').Value = '  +7.30%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.117'
$ws.Range('E43').Value = '  +1.47%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.62'
$ws.Range('E44').Value = '  +8.15%  '
$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '121.17'
$ws.Range('E45').Value = '  +0.58%  '
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.55'
$ws.Range('E47').Value = '  +3.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.198.25'
$ws.Range('E48').Value = '  +4.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.268'
$ws.Range('E49').Value = '  +24.78%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0342'
$ws.Range('E50').Value = '  +11.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.970'
$ws.Range('E51').Value = '  +2.08%  '
